# Insert two new rows (they will become the new rows 424 and 425), pushing the
# existing rows 424-445 down to 426-447.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A424:A425").EntireRow.Insert()

# Populate the first newly inserted row (424)
$ws.Range("A424").Value = 11
$ws.Range("B424").Value = "Vega Monumental Concepción"
$ws.Range("C424").Value = "Bíobío"
$ws.Range("D424").Value = 45147
$ws.Range("E424").Value = 8
$ws.Range("F424").Value = 100114013
$ws.Range("G424").Value = "Zanahoria"
$ws.Range("H424").Value = "Sin especificar"
$ws.Range("I424").Value = "Primera"
$ws.Range("J424").Value = 400
$ws.Range("K424").Value = 4500
$ws.Range("L424").Value = 5000
$ws.Range("M424").Value = 4750
$ws.Range("N424").Value = "$/saco 20 kilos"
$ws.Range("O424").Value = "Región de Ñuble"
$ws.Range("P424").Value = 238
$ws.Range("Q424").Value = 20
$ws.Range("R424").Value = "Hortaliza"

# Populate the second newly inserted row (425)
$ws.Range("A425").Value = 11
$ws.Range("B425").Value = "Vega Monumental Concepción"
$ws.Range("C425").Value = "Bíobío"
$ws.Range("D425").Value = 45147
$ws.Range("E425").Value = 8
$ws.Range("F425").Value = 100114013
$ws.Range("G425").Value = "Zanahoria"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Segunda"
$ws.Range("J425").Value = 200
$ws.Range("K425").Value = 4000
$ws.Range("L425").Value = 4000
$ws.Range("M425").Value = 4000
$ws.Range("N425").Value = "$/saco 20 kilos"
$ws.Range("O425").Value = "Región de Ñuble"
$ws.Range("P425").Value = 200
$ws.Range("Q425").Value = 20
$ws.Range("R425").Value = "Hortaliza"

$wb.Save()
